$wb = $excel.ActiveWorkbook

# Rename sheet "sample3" to "any name you want"
$ws3 = $wb.Worksheets.Item("sample3")
$ws3.Name = "any name you want"

# Activate the renamed sheet (it becomes the selected/active tab)
$ws3.Activate()

# Set scroll position on "sample5" (topLeftCell A82)
$ws5 = $wb.Worksheets.Item("sample5")
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1

# Re-activate "any name you want" sheet to be the final active tab
$ws3.Activate()
